# Applies the figure5 ROI slide revisions:
#  1) Merge the multi-run "(f) DSI-11-Gmax225 In Vivo" caption into a single run.
#  2) Merge the multi-run "(g) DSI-11-Gmax40 In Vivo" caption into a single run.
#  3) Ungroup the outer "组合 15" wrapper group so its three children
#     (the ROI overlay group, the screenshot picture, and the "a" label
#     textbox) become direct top-level shapes on the slide.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- 1 & 2: collapse the split text runs into single runs ------------------
# The host preserves existing runs untouched when the assigned text already
# matches the concatenation of the existing runs, so we first set the text to
# a throwaway placeholder (forcing the whole paragraph to be rewritten) and
# then assign the real caption text, which lands in one run with the
# original run formatting.

$grp20 = $s.Shapes.Item("组合 20")
$capF = $grp20.GroupItems.Item("TextBox 72")
$capF.TextFrame.TextRange.Text = "_"
$capF.TextFrame.TextRange.Text = "(f) DSI-11-Gmax225 In Vivo"

$grp21 = $s.Shapes.Item("组合 21")
$capG = $grp21.GroupItems.Item("TextBox 72")
$capG.TextFrame.TextRange.Text = "_"
$capG.TextFrame.TextRange.Text = "(g) DSI-11-Gmax40 In Vivo"

# --- 3: ungroup "组合 15" ----------------------------------------------------
$wrapper = $s.Shapes.Item("组合 15")
$wrapper.Ungroup() | Out-Null
